$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03644533333333334
$ws.Range("H2").Value = 0.109336
$ws.Range("I2").Value = 0.005561955322140003
$ws.Range("J2").Value = 0.005561955322140003
$ws.Range("M2").Value = 0.3613943333333333
$ws.Range("N2").Value = 1.084183
$ws.Range("O2").Value = 0.194831931013823
$ws.Range("P2").Value = 0.194831931013823
$ws.Range("Q2").Value = 0.01317113694311111
$ws.Range("R2").Value = 0.118540232488
$ws.Range("S2").Value = 0.001083646495625147
$ws.Range("T2").Value = 0.001083646495625147
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03644533333333334
$ws.Range("H3").Value = 0.109336
$ws.Range("I3").Value = 0.005561955322140003
$ws.Range("J3").Value = 0.005561955322140003
$ws.Range("O3").Value = 0.5804398756520781
$ws.Range("P3").Value = 0.5804398756520782
$ws.Range("Q3").Value = 0.03923922043822223
$ws.Range("R3").Value = 0.353152983944
$ws.Range("S3").Value = 0.003228380655565357
$ws.Range("T3").Value = 0.003228380655565358
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03644533333333334
$ws.Range("H4").Value = 0.109336
$ws.Range("I4").Value = 0.005561955322140003
$ws.Range("J4").Value = 0.005561955322140003
$ws.Range("M4").Value = 0.416849
$ws.Range("N4").Value = 1.250547
$ws.Range("O4").Value = 0.2247281933340989
$ws.Range("P4").Value = 0.2247281933340989
$ws.Range("Q4").Value = 0.01519220075466667
$ws.Range("R4").Value = 0.136729806792
$ws.Range("S4").Value = 0.001249928170949498
$ws.Range("T4").Value = 0.001249928170949499
$ws.Range("G5").Value = 6.516166333333334
$ws.Range("H5").Value = 19.548499
$ws.Range("I5").Value = 0.99443804467786
$ws.Range("J5").Value = 0.99443804467786
$ws.Range("M5").Value = 0.3613943333333333
$ws.Range("N5").Value = 1.084183
$ws.Range("O5").Value = 0.194831931013823
$ws.Range("P5").Value = 0.194831931013823
$ws.Range("Q5").Value = 2.354905587924111
$ws.Range("R5").Value = 21.194150291317
$ws.Range("S5").Value = 0.1937482845181978
$ws.Range("T5").Value = 0.1937482845181978
$ws.Range("G6").Value = 6.516166333333334
$ws.Range("H6").Value = 19.548499
$ws.Range("I6").Value = 0.99443804467786
$ws.Range("J6").Value = 0.99443804467786
$ws.Range("O6").Value = 0.5804398756520781
$ws.Range("P6").Value = 0.5804398756520782
$ws.Range("Q6").Value = 7.015693472391224
$ws.Range("R6").Value = 63.141241251521
$ws.Range("S6").Value = 0.5772114949965128
$ws.Range("T6").Value = 0.5772114949965129
$ws.Range("G7").Value = 6.516166333333334
$ws.Range("H7").Value = 19.548499
$ws.Range("I7").Value = 0.99443804467786
$ws.Range("J7").Value = 0.99443804467786
$ws.Range("M7").Value = 0.416849
$ws.Range("N7").Value = 1.250547
$ws.Range("O7").Value = 0.2247281933340989
$ws.Range("P7").Value = 0.2247281933340989
$ws.Range("Q7").Value = 2.716257419883667
$ws.Range("R7").Value = 24.446316778953
$ws.Range("S7").Value = 0.2234782651631494
$ws.Range("T7").Value = 0.2234782651631494
